# ---------------------------------------------------------------------------
# romania_liga-2_2023-2024 update script (commit: "Atualizado por script em
# 26-11-2023 20:30")
#
# The source diff shows two kinds of changes to Sheet1:
#   1) A handful of existing rows had their match-data columns (F:V — home
#      team through match URL; columns A:E — Indice/pais/torneio/temporada/
#      data_partida — are untouched) permuted among themselves. This happens
#      in five independent groups of rows:
#         {2,4}                       simple swap
#         {73,74,75,76,77}            rotation
#         {83,84,85,86,87,88}         rotation (single 6-cycle)
#         {98,99}                     simple swap
#         {103,104,105,107,108,109}   rotation (single 6-cycle; 106 untouched)
#         {112,113,115,116,117}       two rotations: {112,116} and {113,117,115}
#                                     (114 untouched)
#   2) Seven brand-new match rows (134-140, Indice 133-139) are appended,
#      growing the used range from A1:V133 to A1:V140.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the "before" F:V contents of every row that is part of
# --- a permutation, before any of them gets overwritten.
$v2   = $ws.Range("F2:V2").Value()
$v4   = $ws.Range("F4:V4").Value()

$v73  = $ws.Range("F73:V73").Value()
$v74  = $ws.Range("F74:V74").Value()
$v75  = $ws.Range("F75:V75").Value()
$v76  = $ws.Range("F76:V76").Value()
$v77  = $ws.Range("F77:V77").Value()

$v83  = $ws.Range("F83:V83").Value()
$v84  = $ws.Range("F84:V84").Value()
$v85  = $ws.Range("F85:V85").Value()
$v86  = $ws.Range("F86:V86").Value()
$v87  = $ws.Range("F87:V87").Value()
$v88  = $ws.Range("F88:V88").Value()

$v98  = $ws.Range("F98:V98").Value()
$v99  = $ws.Range("F99:V99").Value()

$v103 = $ws.Range("F103:V103").Value()
$v104 = $ws.Range("F104:V104").Value()
$v105 = $ws.Range("F105:V105").Value()
$v107 = $ws.Range("F107:V107").Value()
$v108 = $ws.Range("F108:V108").Value()
$v109 = $ws.Range("F109:V109").Value()

$v112 = $ws.Range("F112:V112").Value()
$v113 = $ws.Range("F113:V113").Value()
$v115 = $ws.Range("F115:V115").Value()
$v116 = $ws.Range("F116:V116").Value()
$v117 = $ws.Range("F117:V117").Value()

# --- Step 2: write each row's new F:V content from the snapshot above.
$ws.Range("F2:V2").Value   = $v4
$ws.Range("F4:V4").Value   = $v2

$ws.Range("F73:V73").Value = $v77
$ws.Range("F74:V74").Value = $v73
$ws.Range("F75:V75").Value = $v74
$ws.Range("F76:V76").Value = $v75
$ws.Range("F77:V77").Value = $v76

$ws.Range("F83:V83").Value = $v87
$ws.Range("F84:V84").Value = $v86
$ws.Range("F85:V85").Value = $v88
$ws.Range("F86:V86").Value = $v85
$ws.Range("F87:V87").Value = $v84
$ws.Range("F88:V88").Value = $v83

$ws.Range("F98:V98").Value = $v99
$ws.Range("F99:V99").Value = $v98

$ws.Range("F103:V103").Value = $v109
$ws.Range("F104:V104").Value = $v108
$ws.Range("F105:V105").Value = $v107
$ws.Range("F107:V107").Value = $v104
$ws.Range("F108:V108").Value = $v103
$ws.Range("F109:V109").Value = $v105

$ws.Range("F112:V112").Value = $v116
$ws.Range("F113:V113").Value = $v117
$ws.Range("F115:V115").Value = $v113
$ws.Range("F116:V116").Value = $v112
$ws.Range("F117:V117").Value = $v115

# --- Step 3: append the 7 new fixture rows (134-140). First clone the
# --- formatting of the last existing data row (133: bold/boxed Indice cell,
# --- date-formatted data_partida cell) down across the new range, then fill
# --- in values cell by cell.
$ws.Range("A133:V133").Copy()
$ws.Range("A134:V140").PasteSpecial(-4122)

# Row 134 (Indice=133)
$ws.Cells.Item(134,1).Value = 133
$ws.Cells.Item(134,2).Value = "romania"
$ws.Cells.Item(134,3).Value = "liga-2"
$ws.Cells.Item(134,4).Value = "2023-2024"
$ws.Cells.Item(134,5).Value = 45255.41666666666
$ws.Cells.Item(134,6).Value = "Mioveni"
$ws.Cells.Item(134,7).Value = 0
$ws.Cells.Item(134,8).Value = "Metaloglobus Bucharest"
$ws.Cells.Item(134,9).Value = 2
$ws.Cells.Item(134,10).Value = 1.86
$ws.Cells.Item(134,11).Value = "23/11/2023 22:12"
$ws.Cells.Item(134,12).Value = 2.02
$ws.Cells.Item(134,13).Value = "25/11/2023 09:51"
$ws.Cells.Item(134,14).Value = 3.22
$ws.Cells.Item(134,15).Value = "23/11/2023 22:12"
$ws.Cells.Item(134,16).Value = 3.25
$ws.Cells.Item(134,17).Value = "25/11/2023 09:51"
$ws.Cells.Item(134,18).Value = 3.93
$ws.Cells.Item(134,19).Value = "23/11/2023 22:12"
$ws.Cells.Item(134,20).Value = 3.91
$ws.Cells.Item(134,21).Value = "25/11/2023 09:51"
$ws.Cells.Item(134,22).Value = "https://www.betexplorer.com/football/romania/liga-2/mioveni-metaloglobus-bucharest/CWDoxORd/"

# Row 135 (Indice=134)
$ws.Cells.Item(135,1).Value = 134
$ws.Cells.Item(135,2).Value = "romania"
$ws.Cells.Item(135,3).Value = "liga-2"
$ws.Cells.Item(135,4).Value = "2023-2024"
$ws.Cells.Item(135,5).Value = 45255.41666666666
$ws.Cells.Item(135,6).Value = "Csikszereda M. Ciuc"
$ws.Cells.Item(135,7).Value = 1
$ws.Cells.Item(135,8).Value = "CSM Slatina"
$ws.Cells.Item(135,9).Value = 0
$ws.Cells.Item(135,10).Value = 1.94
$ws.Cells.Item(135,11).Value = "25/11/2023 00:42"
$ws.Cells.Item(135,12).Value = 2.11
$ws.Cells.Item(135,13).Value = "25/11/2023 09:58"
$ws.Cells.Item(135,14).Value = 3.31
$ws.Cells.Item(135,15).Value = "25/11/2023 00:42"
$ws.Cells.Item(135,16).Value = 3.13
$ws.Cells.Item(135,17).Value = "25/11/2023 09:58"
$ws.Cells.Item(135,18).Value = 4.04
$ws.Cells.Item(135,19).Value = "25/11/2023 00:42"
$ws.Cells.Item(135,20).Value = 3.81
$ws.Cells.Item(135,21).Value = "25/11/2023 09:58"
$ws.Cells.Item(135,22).Value = "https://www.betexplorer.com/football/romania/liga-2/miercurea-ciuc-csm-slatina/ITfOKMdM/"

# Row 136 (Indice=135)
$ws.Cells.Item(136,1).Value = 135
$ws.Cells.Item(136,2).Value = "romania"
$ws.Cells.Item(136,3).Value = "liga-2"
$ws.Cells.Item(136,4).Value = "2023-2024"
$ws.Cells.Item(136,5).Value = 45255.41666666666
$ws.Cells.Item(136,6).Value = "Progresul Spartac"
$ws.Cells.Item(136,7).Value = 1
$ws.Cells.Item(136,8).Value = "Steaua Bucuresti"
$ws.Cells.Item(136,9).Value = 6
$ws.Cells.Item(136,10).Value = 7.1
$ws.Cells.Item(136,11).Value = "23/11/2023 22:12"
$ws.Cells.Item(136,12).Value = 11.17
$ws.Cells.Item(136,13).Value = "25/11/2023 09:51"
$ws.Cells.Item(136,14).Value = 4.56
$ws.Cells.Item(136,15).Value = "23/11/2023 22:12"
$ws.Cells.Item(136,16).Value = 5.61
$ws.Cells.Item(136,17).Value = "25/11/2023 09:51"
$ws.Cells.Item(136,18).Value = 1.37
$ws.Cells.Item(136,19).Value = "23/11/2023 22:12"
$ws.Cells.Item(136,20).Value = 1.27
$ws.Cells.Item(136,21).Value = "25/11/2023 09:51"
$ws.Cells.Item(136,22).Value = "https://www.betexplorer.com/football/romania/liga-2/progresul-spartac-csa-steaua-bucuresti/0hpTJtCS/"

# Row 137 (Indice=136)
$ws.Cells.Item(137,1).Value = 136
$ws.Cells.Item(137,2).Value = "romania"
$ws.Cells.Item(137,3).Value = "liga-2"
$ws.Cells.Item(137,4).Value = "2023-2024"
$ws.Cells.Item(137,5).Value = 45255.41666666666
$ws.Cells.Item(137,6).Value = "CSM Resita"
$ws.Cells.Item(137,7).Value = 3
$ws.Cells.Item(137,8).Value = "Alexandria"
$ws.Cells.Item(137,9).Value = 1
$ws.Cells.Item(137,10).Value = 1.68
$ws.Cells.Item(137,11).Value = "23/11/2023 22:12"
$ws.Cells.Item(137,12).Value = 1.63
$ws.Cells.Item(137,13).Value = "25/11/2023 09:55"
$ws.Cells.Item(137,14).Value = 3.61
$ws.Cells.Item(137,15).Value = "23/11/2023 22:12"
$ws.Cells.Item(137,16).Value = 3.87
$ws.Cells.Item(137,17).Value = "25/11/2023 09:58"
$ws.Cells.Item(137,18).Value = 4.34
$ws.Cells.Item(137,19).Value = "23/11/2023 22:12"
$ws.Cells.Item(137,20).Value = 5.35
$ws.Cells.Item(137,21).Value = "25/11/2023 09:57"
$ws.Cells.Item(137,22).Value = "https://www.betexplorer.com/football/romania/liga-2/csm-resita-csm-alexandria/WzyaEKZk/"

# Row 138 (Indice=137)
$ws.Cells.Item(138,1).Value = 137
$ws.Cells.Item(138,2).Value = "romania"
$ws.Cells.Item(138,3).Value = "liga-2"
$ws.Cells.Item(138,4).Value = "2023-2024"
$ws.Cells.Item(138,5).Value = 45255.4375
$ws.Cells.Item(138,6).Value = "Tunari"
$ws.Cells.Item(138,7).Value = 0
$ws.Cells.Item(138,8).Value = "Chindia Targoviste"
$ws.Cells.Item(138,9).Value = 2
$ws.Cells.Item(138,10).Value = 4.21
$ws.Cells.Item(138,11).Value = "25/11/2023 01:42"
$ws.Cells.Item(138,12).Value = 4.49
$ws.Cells.Item(138,13).Value = "25/11/2023 10:27"
$ws.Cells.Item(138,14).Value = 3.45
$ws.Cells.Item(138,15).Value = "25/11/2023 01:42"
$ws.Cells.Item(138,16).Value = 3.45
$ws.Cells.Item(138,17).Value = "25/11/2023 10:27"
$ws.Cells.Item(138,18).Value = 1.83
$ws.Cells.Item(138,19).Value = "25/11/2023 01:42"
$ws.Cells.Item(138,20).Value = 1.83
$ws.Cells.Item(138,21).Value = "25/11/2023 10:27"
$ws.Cells.Item(138,22).Value = "https://www.betexplorer.com/football/romania/liga-2/tunari-chindia-targoviste/AHWlyrs3/"

# Row 139 (Indice=138)
$ws.Cells.Item(139,1).Value = 138
$ws.Cells.Item(139,2).Value = "romania"
$ws.Cells.Item(139,3).Value = "liga-2"
$ws.Cells.Item(139,4).Value = "2023-2024"
$ws.Cells.Item(139,5).Value = 45256.41666666666
$ws.Cells.Item(139,6).Value = "Viitorul Tg. Jiu"
$ws.Cells.Item(139,7).Value = 2
$ws.Cells.Item(139,8).Value = "Unirea Dej"
$ws.Cells.Item(139,9).Value = 1
$ws.Cells.Item(139,10).Value = 2.16
$ws.Cells.Item(139,11).Value = "24/11/2023 22:12"
$ws.Cells.Item(139,12).Value = 2.08
$ws.Cells.Item(139,13).Value = "26/11/2023 09:56"
$ws.Cells.Item(139,14).Value = 2.98
$ws.Cells.Item(139,15).Value = "24/11/2023 22:12"
$ws.Cells.Item(139,16).Value = 3.18
$ws.Cells.Item(139,17).Value = "26/11/2023 09:57"
$ws.Cells.Item(139,18).Value = 3.29
$ws.Cells.Item(139,19).Value = "24/11/2023 22:12"
$ws.Cells.Item(139,20).Value = 3.83
$ws.Cells.Item(139,21).Value = "26/11/2023 09:56"
$ws.Cells.Item(139,22).Value = "https://www.betexplorer.com/football/romania/liga-2/viitorul-targu-jiu-unirea-dej/ARZ2Dvkd/"

# Row 140 (Indice=139)
$ws.Cells.Item(140,1).Value = 139
$ws.Cells.Item(140,2).Value = "romania"
$ws.Cells.Item(140,3).Value = "liga-2"
$ws.Cells.Item(140,4).Value = "2023-2024"
$ws.Cells.Item(140,5).Value = 45256.45833333334
$ws.Cells.Item(140,6).Value = "Hunedoara"
$ws.Cells.Item(140,7).Value = 0
$ws.Cells.Item(140,8).Value = "Selimbar"
$ws.Cells.Item(140,9).Value = 0
$ws.Cells.Item(140,10).Value = 1.69
$ws.Cells.Item(140,11).Value = "24/11/2023 23:13"
$ws.Cells.Item(140,12).Value = 1.77
$ws.Cells.Item(140,13).Value = "26/11/2023 10:42"
$ws.Cells.Item(140,14).Value = 3.43
$ws.Cells.Item(140,15).Value = "24/11/2023 23:13"
$ws.Cells.Item(140,16).Value = 3.5
$ws.Cells.Item(140,17).Value = "26/11/2023 10:53"
$ws.Cells.Item(140,18).Value = 4.51
$ws.Cells.Item(140,19).Value = "24/11/2023 23:13"
$ws.Cells.Item(140,20).Value = 4.85
$ws.Cells.Item(140,21).Value = "26/11/2023 10:42"
$ws.Cells.Item(140,22).Value = "https://www.betexplorer.com/football/romania/liga-2/corvinul-hunedoara-selimbar/21WfF0Kq/"
